# The workbook contains a daily price log for "Piña" (pineapple) at
# "Vega Monumental Concepción". A new price record (dated 2022-01-20,
# serial 44581) is inserted as row 40, pushing the existing rows 40-135
# down to 41-136 (dimension grows from A1:T135 to A1:T136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 40 - this shifts every row that was
# at 40..135 down to 41..136, preserving all of their data/formatting.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new daily record.
$ws.Cells.Item(40, 1).Value  = 11
$ws.Cells.Item(40, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value  = "Bíobío"
$ws.Cells.Item(40, 4).Value  = 44581
$ws.Cells.Item(40, 5).Value  = 8
$ws.Cells.Item(40, 6).Value  = "Fruta"
$ws.Cells.Item(40, 7).Value  = 100108
$ws.Cells.Item(40, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(40, 9).Value  = 100108005
$ws.Cells.Item(40, 10).Value = "Piña"
$ws.Cells.Item(40, 11).Value = "Caramelo"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 310
$ws.Cells.Item(40, 14).Value = 12000
$ws.Cells.Item(40, 15).Value = 13000
$ws.Cells.Item(40, 16).Value = 12484
$ws.Cells.Item(40, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(40, 18).Value = "Ecuador"
$ws.Cells.Item(40, 19).Value = 892
$ws.Cells.Item(40, 20).Value = 14

# Match the date cell formatting/style used by the other date cells in
# column D (style index 2 in the original sheet, a date number format).
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
